$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 4; this pushes the existing rows 4-8 down to 5-9
$ws.Rows.Item(4).Insert()

# Populate the newly inserted row 4 with the new weekly record
$ws.Cells.Item(4, 1).Value = 7
$ws.Cells.Item(4, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(4, 3).Value = "Ñuble"
$ws.Cells.Item(4, 4).Value = 44645
$ws.Cells.Item(4, 4).NumberFormat = $ws.Cells.Item(5, 4).NumberFormat
$ws.Cells.Item(4, 5).Value = 16
$ws.Cells.Item(4, 6).Value = 100112043
$ws.Cells.Item(4, 7).Value = "Pepino dulce"
$ws.Cells.Item(4, 8).Value = "Cultivar IV Región"
$ws.Cells.Item(4, 9).Value = "Primera"
$ws.Cells.Item(4, 10).Value = 60
$ws.Cells.Item(4, 11).Value = 15000
$ws.Cells.Item(4, 12).Value = 16000
$ws.Cells.Item(4, 13).Value = 15500
$ws.Cells.Item(4, 14).Value = "`$/bandeja 18 kilos"
$ws.Cells.Item(4, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(4, 16).Value = 861
$ws.Cells.Item(4, 17).Value = 18
$ws.Cells.Item(4, 18).Value = "Hortaliza"
